# TC42_Canine_Filter_Breed-ShihTzu.xlsx — "corrected ICDC Breed 1-14 scripts"
#
# The FilesTab Cypher query (cell B4 on the "startup" sheet) is corrected:
#   - the `coalesce(f.file_type, '') AS `File Type`` return column is removed
#   - the `coalesce(demo.breed,'') AS Breed ,` return column is removed
#   - the following `Diagnosis` line keeps its content but gains one extra
#     leading space of indentation (artifact of the two preceding lines
#     being deleted around it)
# The active-cell selection on that sheet moves from C13 to B4, and the
# row height for row 4 shrinks to match its new (shorter) wrapped content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Shih Tzu'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesQuery

# Row 4 used wrapped text that was 246.5pt tall for the old (longer) query;
# the trimmed query only needs 217.5pt — the same height row 3 uses for its
# similarly-sized query text.
$ws.Rows.Item(4).RowHeight = 217.5

# Move the sheet's saved selection/active cell from C13 to B4.
$ws.Activate() | Out-Null
$ws.Range("B4").Select() | Out-Null
